$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Margin" column (column E) entirely.
$ws.Range("E1:E6").EntireColumn.Delete()

# The item/price columns store their numbers as plain text (shared strings),
# not real numbers. Temporarily format the price columns as Text so the
# new values are not auto-converted to numbers, then restore the "Normal"
# style so no extra cell formatting lingers on the cells.
$priceRange = $ws.Range("B2:D6")
$priceRange.NumberFormat = "@"

# Update the price values (Buy/Sell/Current all share the same value per row).
$ws.Range("B2").Value = "236"
$ws.Range("C2").Value = "236"
$ws.Range("D2").Value = "236"

$ws.Range("B3").Value = "323"
$ws.Range("C3").Value = "323"
$ws.Range("D3").Value = "323"

$ws.Range("B4").Value = "2,053"
$ws.Range("C4").Value = "2,053"
$ws.Range("D4").Value = "2,053"

$ws.Range("B5").Value = "87"
$ws.Range("C5").Value = "87"
$ws.Range("D5").Value = "87"

$ws.Range("B6").Value = "174"
$ws.Range("C6").Value = "174"
$ws.Range("D6").Value = "174"

# Drop the temporary Text formatting now that the values are locked in as
# text, so the cells end up back on the plain default style.
$priceRange.Style = "Normal"
